$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4796.385
$ws.Range("J17").Value = 1862.75
$ws.Range("L17").Value = 5588.25
$ws.Range("N17").Value = -5924.25

$ws.Range("H19").Value = 3705.4092
$ws.Range("I19").Value = 919.9167
$ws.Range("K19").Value = 919.9167
$ws.Range("M19").Value = -744.9167

$ws.Range("H105").Value = 76389.836
$ws.Range("J105").Value = 77668
$ws.Range("L105").Value = 77668
$ws.Range("N105").Value = -84656

$ws.Range("H112").Value = 1559.9
$ws.Range("J112").Value = 1559.9
$ws.Range("L112").Value = 4679.700000000001
$ws.Range("N112").Value = -6895.700000000001

$ws.Range("H116").Value = 97718.086
$ws.Range("I116").Value = 186969.83
$ws.Range("J116").Value = 8466.333000000001
$ws.Range("K116").Value = 186969.83
$ws.Range("L116").Value = 8466.333000000001
$ws.Range("M116").Value = -183527.83
$ws.Range("N116").Value = -15350.333

$ws.Range("H132").Value = 2375.1406
$ws.Range("I132").Value = 1722.4706
$ws.Range("K132").Value = 5167.4118
$ws.Range("M132").Value = -2637.4118

$ws.Range("H138").Value = 1663.0377
$ws.Range("I138").Value = 1249.8
$ws.Range("J138").Value = 1913.4849
$ws.Range("K138").Value = 3749.4
$ws.Range("L138").Value = 5740.4547
$ws.Range("M138").Value = 1390.6
$ws.Range("N138").Value = -16020.4547

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 193
$ws.Range("I5").Value = 193
$ws.Range("K5").Value = 193
$ws.Range("M5").Value = -81

$ws.Range("H14").Value = 457
$ws.Range("I14").Value = 299.8
$ws.Range("K14").Value = 299.8
$ws.Range("M14").Value = -124.8

$ws.Range("H61").Value = 2334.5881
$ws.Range("I61").Value = 2012.5333
$ws.Range("J61").Value = 4750
$ws.Range("K61").Value = 2012.5333
$ws.Range("L61").Value = 4750
$ws.Range("M61").Value = -1800.5333
$ws.Range("N61").Value = -5174

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H112").Value = 37974.5
$ws.Range("J112").Value = 37974.5
$ws.Range("L112").Value = 37974.5
$ws.Range("N112").Value = -40928.5

$ws.Range("H136").Value = 2334.5881
$ws.Range("I136").Value = 2012.5333
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 6037.5999
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -3487.5999
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 193
$ws.Range("I4").Value = 193
$ws.Range("K4").Value = 193
$ws.Range("M4").Value = -78

$ws.Range("H15").Value = 450
$ws.Range("J15").Value = 450
$ws.Range("L15").Value = 450
$ws.Range("N15").Value = -904

$ws.Range("H76").Value = 16078.5
$ws.Range("J76").Value = 16078.5
$ws.Range("L76").Value = 16078.5
$ws.Range("N76").Value = -16708.5

$ws.Range("H79").Value = 16078.5
$ws.Range("J79").Value = 16078.5
$ws.Range("L79").Value = 16078.5
$ws.Range("N79").Value = -18262.5

$ws.Range("H86").Value = 2421.4285
$ws.Range("J86").Value = 3129.923
$ws.Range("L86").Value = 3129.923
$ws.Range("N86").Value = -5375.923

$ws.Range("H88").Value = 18296.143
$ws.Range("J88").Value = 18296.143
$ws.Range("L88").Value = 18296.143
$ws.Range("N88").Value = -19108.143

$ws.Range("H89").Value = 2421.4285
$ws.Range("J89").Value = 3129.923
$ws.Range("L89").Value = 15649.615
$ws.Range("N89").Value = -26881.615

$ws.Range("H91").Value = 18296.143
$ws.Range("J91").Value = 18296.143
$ws.Range("L91").Value = 18296.143
$ws.Range("N91").Value = -21104.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2897.5317
$ws.Range("I31").Value = 1520.375
$ws.Range("K31").Value = 1520.375
$ws.Range("M31").Value = -1225.375

$ws.Range("H34").Value = 2897.5317
$ws.Range("I34").Value = 1520.375
$ws.Range("K34").Value = 1520.375
$ws.Range("M34").Value = -1318.375

$ws.Range("H62").Value = 79056.875
$ws.Range("J62").Value = 89636.42999999999
$ws.Range("L62").Value = 89636.42999999999
$ws.Range("N62").Value = -90884.42999999999

$ws.Range("H65").Value = 79056.875
$ws.Range("J65").Value = 89636.42999999999
$ws.Range("L65").Value = 448182.15
$ws.Range("N65").Value = -454422.15

$ws.Range("H88").Value = 15050
$ws.Range("J88").Value = 15050
$ws.Range("L88").Value = 15050
$ws.Range("N88").Value = -15862

$ws.Range("H91").Value = 15050
$ws.Range("J91").Value = 15050
$ws.Range("L91").Value = 15050
$ws.Range("N91").Value = -17858

$ws.Range("H117").Value = 100000
$ws.Range("J117").Value = 100000
$ws.Range("L117").Value = 100000
$ws.Range("N117").Value = -109178

$ws.Range("H132").Value = 2879.2307
$ws.Range("I132").Value = 1867.3636
$ws.Range("J132").Value = 8444.5
$ws.Range("K132").Value = 5602.0908
$ws.Range("L132").Value = 25333.5
$ws.Range("M132").Value = -3072.0908
$ws.Range("N132").Value = -30393.5

$ws.Range("H134").Value = 3238.258
$ws.Range("I134").Value = 2680.7036
$ws.Range("K134").Value = 8042.110799999999
$ws.Range("M134").Value = -5507.110799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 6459.25
$ws.Range("J81").Value = 8127
$ws.Range("L81").Value = 24381
$ws.Range("N81").Value = -26627

$ws.Range("H84").Value = 6459.25
$ws.Range("J84").Value = 8127
$ws.Range("L84").Value = 73143
$ws.Range("N84").Value = -84375

$ws.Range("H106").Value = 5999.875
$ws.Range("I106").Value = 2999.5
$ws.Range("J106").Value = 7000
$ws.Range("K106").Value = 8998.5
$ws.Range("L106").Value = 21000
$ws.Range("M106").Value = -8052.5
$ws.Range("N106").Value = -22892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6838.0586
$ws.Range("I70").Value = 6183.2
$ws.Range("K70").Value = 6183.2
$ws.Range("M70").Value = -5913.2

$ws.Range("H73").Value = 6838.0586
$ws.Range("I73").Value = 6183.2
$ws.Range("K73").Value = 6183.2
$ws.Range("M73").Value = -5247.2

$ws.Range("H92").Value = 14348.75
$ws.Range("J92").Value = 14348.75
$ws.Range("L92").Value = 14348.75
$ws.Range("N92").Value = -18092.75

$ws.Range("H126").Value = 6506.6177
$ws.Range("I126").Value = 9411.706
$ws.Range("J126").Value = 3601.5293
$ws.Range("K126").Value = 28235.118
$ws.Range("L126").Value = 10804.5879
$ws.Range("M126").Value = -25765.118
$ws.Range("N126").Value = -15744.5879

$ws.Range("H132").Value = 4304.381
$ws.Range("I132").Value = 5123
$ws.Range("K132").Value = 15369
$ws.Range("M132").Value = -12839

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 986.4
$ws.Range("I30").Value = 986.4
$ws.Range("K30").Value = 986.4
$ws.Range("M30").Value = -878.4

$ws.Range("H35").Value = 895.0909
$ws.Range("I35").Value = 884.6
$ws.Range("K35").Value = 884.6
$ws.Range("M35").Value = -548.6

$ws.Range("H61").Value = 5757.5
$ws.Range("I61").Value = 7065.091
$ws.Range("J61").Value = 2880.8
$ws.Range("K61").Value = 7065.091
$ws.Range("L61").Value = 2880.8
$ws.Range("M61").Value = -6863.091
$ws.Range("N61").Value = -3284.8

$ws.Range("H68").Value = 7400.4
$ws.Range("I68").Value = 6001.5
$ws.Range("K68").Value = 6001.5
$ws.Range("M68").Value = -5252.5

$ws.Range("H71").Value = 7400.4
$ws.Range("I71").Value = 6001.5
$ws.Range("K71").Value = 30007.5
$ws.Range("M71").Value = -26263.5

$ws.Range("H82").Value = 1668.25
$ws.Range("J82").Value = 1766.3334
$ws.Range("L82").Value = 1766.3334
$ws.Range("N82").Value = -2488.3334

$ws.Range("H85").Value = 1668.25
$ws.Range("J85").Value = 1766.3334
$ws.Range("L85").Value = 1766.3334
$ws.Range("N85").Value = -4262.3334

$ws.Range("H93").Value = 15876613
$ws.Range("J93").Value = 5001
$ws.Range("L93").Value = 5001
$ws.Range("N93").Value = -7497

$ws.Range("H113").Value = 5757.5
$ws.Range("I113").Value = 7065.091
$ws.Range("J113").Value = 2880.8
$ws.Range("K113").Value = 7065.091
$ws.Range("L113").Value = 2880.8
$ws.Range("M113").Value = -4895.091
$ws.Range("N113").Value = -7220.8

$ws.Range("H132").Value = 27032146
$ws.Range("I132").Value = 41669890
$ws.Range("J132").Value = 8622.538
$ws.Range("K132").Value = 125009670
$ws.Range("L132").Value = 25867.614
$ws.Range("M132").Value = -125007140
$ws.Range("N132").Value = -30927.614

$ws.Range("H136").Value = 5210.8076
$ws.Range("I136").Value = 6167.778
$ws.Range("K136").Value = 18503.334
$ws.Range("M136").Value = -15953.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 44500
$ws.Range("J21").Value = 44500
$ws.Range("L21").Value = 44500
$ws.Range("N21").Value = -44970

$ws.Range("H35").Value = 44500
$ws.Range("J35").Value = 44500
$ws.Range("L35").Value = 44500
$ws.Range("N35").Value = -45080

$ws.Range("H122").Value = 1818.56
$ws.Range("I122").Value = 1612.619
$ws.Range("K122").Value = 4837.857
$ws.Range("M122").Value = -2387.857
